$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full-row rewrites (rows whose entire B:AC content was swapped with another row) ---
$row62 = New-Object "object[,]" 1,28
$row62[0,0] = 6078944
$row62[0,1] = "Chile Primera Division"
$row62[0,2] = "Chile Primera Division"
$row62[0,3] = [datetime]"2023-03-18T18:00:00"
$row62[0,4] = "OHiggins"
$row62[0,5] = "Coquimbo Unido"
$row62[0,6] = 0
$row62[0,7] = 2
$row62[0,8] = "A"
$row62[0,9] = 2
$row62[0,10] = 3.25
$row62[0,11] = 3.4
$row62[0,12] = 2.15
$row62[0,13] = 3.25
$row62[0,14] = 3.6
$row62[0,15] = -0.25
$row62[0,16] = 1.8
$row62[0,17] = 2
$row62[0,18] = 2.25
$row62[0,19] = 1.85
$row62[0,20] = 1.95
$row62[0,21] = -1
$row62[0,22] = -1
$row62[0,23] = 2.6
$row62[0,24] = -1
$row62[0,25] = 1
$row62[0,26] = -0.5
$row62[0,27] = 0.475
$ws.Range("B62:AC62").Value = $row62

$row63 = New-Object "object[,]" 1,28
$row63[0,0] = 6078868
$row63[0,1] = "Chile Primera Division"
$row63[0,2] = "Chile Primera Division"
$row63[0,3] = [datetime]"2023-03-18T18:00:00"
$row63[0,4] = "Cobresal"
$row63[0,5] = "Colo Colo"
$row63[0,6] = 3
$row63[0,7] = 1
$row63[0,8] = "H"
$row63[0,9] = 2.875
$row63[0,10] = 3.25
$row63[0,11] = 2.25
$row63[0,12] = 3.3
$row63[0,13] = 3.3
$row63[0,14] = 2.25
$row63[0,15] = 0.25
$row63[0,16] = 1.875
$row63[0,17] = 1.925
$row63[0,18] = 2.5
$row63[0,19] = 2
$row63[0,20] = 1.8
$row63[0,21] = 2.3
$row63[0,22] = -1
$row63[0,23] = -1
$row63[0,24] = 0.875
$row63[0,25] = -1
$row63[0,26] = 1
$row63[0,27] = -1
$ws.Range("B63:AC63").Value = $row63

$row82 = New-Object "object[,]" 1,28
$row82[0,0] = 6075788
$row82[0,1] = "Chile Primera Division"
$row82[0,2] = "Chile Primera Division"
$row82[0,3] = [datetime]"2023-04-23T16:00:00"
$row82[0,4] = "Magallanes"
$row82[0,5] = "Deportes Copiapo"
$row82[0,6] = 0
$row82[0,7] = 2
$row82[0,8] = "A"
$row82[0,9] = 1.95
$row82[0,10] = 3.4
$row82[0,11] = 3.5
$row82[0,12] = 2
$row82[0,13] = 3.5
$row82[0,14] = 3.75
$row82[0,15] = -0.5
$row82[0,16] = 2.025
$row82[0,17] = 1.825
$row82[0,18] = 2.5
$row82[0,19] = 1.875
$row82[0,20] = 1.975
$row82[0,21] = -1
$row82[0,22] = -1
$row82[0,23] = 2.75
$row82[0,24] = -1
$row82[0,25] = 0.825
$row82[0,26] = -1
$row82[0,27] = 0.9750000000000001
$ws.Range("B82:AC82").Value = $row82

$row83 = New-Object "object[,]" 1,28
$row83[0,0] = 6078871
$row83[0,1] = "Chile Primera Division"
$row83[0,2] = "Chile Primera Division"
$row83[0,3] = [datetime]"2023-04-23T16:00:00"
$row83[0,4] = "Cobresal"
$row83[0,5] = "Curico Unido"
$row83[0,6] = 2
$row83[0,7] = 0
$row83[0,8] = "H"
$row83[0,9] = 1.75
$row83[0,10] = 3.6
$row83[0,11] = 4.2
$row83[0,12] = 1.444
$row83[0,13] = 4.75
$row83[0,14] = 7
$row83[0,15] = -1.25
$row83[0,16] = 1.95
$row83[0,17] = 1.85
$row83[0,18] = 3
$row83[0,19] = 1.825
$row83[0,20] = 1.975
$row83[0,21] = 0.444
$row83[0,22] = -1
$row83[0,23] = -1
$row83[0,24] = 0.95
$row83[0,25] = -1
$row83[0,26] = -1
$row83[0,27] = 0.9750000000000001
$ws.Range("B83:AC83").Value = $row83

$row94 = New-Object "object[,]" 1,28
$row94[0,0] = 6078878
$row94[0,1] = "Chile Primera Division"
$row94[0,2] = "Chile Primera Division"
$row94[0,3] = [datetime]"2023-05-09T19:00:00"
$row94[0,4] = "Union Espanola"
$row94[0,5] = "Union La Calera"
$row94[0,6] = 3
$row94[0,7] = 2
$row94[0,8] = "H"
$row94[0,9] = 2.2
$row94[0,10] = 3.25
$row94[0,11] = 3.1
$row94[0,12] = 2.55
$row94[0,13] = 3.2
$row94[0,14] = 2.875
$row94[0,15] = 0
$row94[0,16] = 1.8
$row94[0,17] = 2.05
$row94[0,18] = 2.5
$row94[0,19] = 1.925
$row94[0,20] = 1.875
$row94[0,21] = 1.55
$row94[0,22] = -1
$row94[0,23] = -1
$row94[0,24] = 0.8
$row94[0,25] = -1
$row94[0,26] = 0.925
$row94[0,27] = -1
$ws.Range("B94:AC94").Value = $row94

$row95 = New-Object "object[,]" 1,28
$row95[0,0] = 6075792
$row95[0,1] = "Chile Primera Division"
$row95[0,2] = "Chile Primera Division"
$row95[0,3] = [datetime]"2023-05-09T19:00:00"
$row95[0,4] = "Nublense"
$row95[0,5] = "Deportes Copiapo"
$row95[0,6] = 1
$row95[0,7] = 0
$row95[0,8] = "H"
$row95[0,9] = 1.85
$row95[0,10] = 3.5
$row95[0,11] = 3.6
$row95[0,12] = 1.95
$row95[0,13] = 3.5
$row95[0,14] = 4
$row95[0,15] = -0.5
$row95[0,16] = 1.9
$row95[0,17] = 1.95
$row95[0,18] = 2.5
$row95[0,19] = 1.975
$row95[0,20] = 1.875
$row95[0,21] = 0.95
$row95[0,22] = -1
$row95[0,23] = -1
$row95[0,24] = 0.8999999999999999
$row95[0,25] = -1
$row95[0,26] = -1
$row95[0,27] = 0.875
$ws.Range("B95:AC95").Value = $row95

$row97 = New-Object "object[,]" 1,28
$row97[0,0] = 6532619
$row97[0,1] = "Chile Primera Division"
$row97[0,2] = "Chile Primera Division"
$row97[0,3] = [datetime]"2023-05-10T19:00:00"
$row97[0,4] = "Coquimbo Unido"
$row97[0,5] = "Universidad de Chile"
$row97[0,6] = 2
$row97[0,7] = 1
$row97[0,8] = "H"
$row97[0,9] = 2.55
$row97[0,10] = 3.4
$row97[0,11] = 2.5
$row97[0,12] = 2.45
$row97[0,13] = 3.3
$row97[0,14] = 3
$row97[0,15] = -0.25
$row97[0,16] = 2.05
$row97[0,17] = 1.75
$row97[0,18] = 2.25
$row97[0,19] = 2
$row97[0,20] = 1.8
$row97[0,21] = 1.45
$row97[0,22] = -1
$row97[0,23] = -1
$row97[0,24] = 1.05
$row97[0,25] = -1
$row97[0,26] = 1
$row97[0,27] = -1
$ws.Range("B97:AC97").Value = $row97

$row98 = New-Object "object[,]" 1,28
$row98[0,0] = 6078876
$row98[0,1] = "Chile Primera Division"
$row98[0,2] = "Chile Primera Division"
$row98[0,3] = [datetime]"2023-05-10T19:00:00"
$row98[0,4] = "Cobresal"
$row98[0,5] = "OHiggins"
$row98[0,6] = 1
$row98[0,7] = 0
$row98[0,8] = "H"
$row98[0,9] = 1.85
$row98[0,10] = 3.6
$row98[0,11] = 3.75
$row98[0,12] = 1.909
$row98[0,13] = 3.6
$row98[0,14] = 4.2
$row98[0,15] = -0.5
$row98[0,16] = 1.85
$row98[0,17] = 1.95
$row98[0,18] = 2.5
$row98[0,19] = 1.875
$row98[0,20] = 1.925
$row98[0,21] = 0.909
$row98[0,22] = -1
$row98[0,23] = -1
$row98[0,24] = 0.8500000000000001
$row98[0,25] = -1
$row98[0,26] = -1
$row98[0,27] = 0.925
$ws.Range("B98:AC98").Value = $row98

$row212 = New-Object "object[,]" 1,28
$row212[0,0] = 7323253
$row212[0,1] = "Chile Primera Division"
$row212[0,2] = "Chile Primera Division"
$row212[0,3] = [datetime]"2023-11-12T20:00:00"
$row212[0,4] = "Union Espanola"
$row212[0,5] = "OHiggins"
$row212[0,6] = 3
$row212[0,7] = 3
$row212[0,8] = "D"
$row212[0,9] = 2
$row212[0,10] = 3.4
$row212[0,11] = 3.5
$row212[0,12] = 2.1
$row212[0,13] = 3.5
$row212[0,14] = 3.75
$row212[0,15] = -0.5
$row212[0,16] = 2.025
$row212[0,17] = 1.775
$row212[0,18] = 2.5
$row212[0,19] = 1.95
$row212[0,20] = 1.85
$row212[0,21] = -1
$row212[0,22] = 2.5
$row212[0,23] = -1
$row212[0,24] = -1
$row212[0,25] = 0.7749999999999999
$row212[0,26] = 0.95
$row212[0,27] = -1
$ws.Range("B212:AC212").Value = $row212

$row213 = New-Object "object[,]" 1,28
$row213[0,0] = 7323186
$row213[0,1] = "Chile Primera Division"
$row213[0,2] = "Chile Primera Division"
$row213[0,3] = [datetime]"2023-11-12T20:00:00"
$row213[0,4] = "Coquimbo Unido"
$row213[0,5] = "Deportes Copiapo"
$row213[0,6] = 1
$row213[0,7] = 0
$row213[0,8] = "H"
$row213[0,9] = 2
$row213[0,10] = 3.4
$row213[0,11] = 3.5
$row213[0,12] = 1.727
$row213[0,13] = 3.8
$row213[0,14] = 4.75
$row213[0,15] = -0.75
$row213[0,16] = 1.9
$row213[0,17] = 1.9
$row213[0,18] = 2.75
$row213[0,19] = 1.85
$row213[0,20] = 1.95
$row213[0,21] = 0.7270000000000001
$row213[0,22] = -1
$row213[0,23] = -1
$row213[0,24] = 0.45
$row213[0,25] = -0.5
$row213[0,26] = -1
$row213[0,27] = 0.95
$ws.Range("B213:AC213").Value = $row213

$row217 = New-Object "object[,]" 1,28
$row217[0,0] = 7494647
$row217[0,1] = "Chile Primera Division"
$row217[0,2] = "Chile Primera Division"
$row217[0,3] = [datetime]"2023-11-25T18:00:00"
$row217[0,4] = "Huachipato"
$row217[0,5] = "Universidad Catolica"
$row217[0,6] = 1
$row217[0,7] = 1
$row217[0,8] = "D"
$row217[0,9] = 2.2
$row217[0,10] = 3.4
$row217[0,11] = 3.2
$row217[0,12] = 1.8
$row217[0,13] = 3.6
$row217[0,14] = 4.333
$row217[0,15] = -0.75
$row217[0,16] = 1.975
$row217[0,17] = 1.875
$row217[0,18] = 2.75
$row217[0,19] = 1.975
$row217[0,20] = 1.875
$row217[0,21] = -1
$row217[0,22] = 2.6
$row217[0,23] = -1
$row217[0,24] = -1
$row217[0,25] = 0.875
$row217[0,26] = -1
$row217[0,27] = 0.875
$ws.Range("B217:AC217").Value = $row217

$row218 = New-Object "object[,]" 1,28
$row218[0,0] = 7494646
$row218[0,1] = "Chile Primera Division"
$row218[0,2] = "Chile Primera Division"
$row218[0,3] = [datetime]"2023-11-25T18:00:00"
$row218[0,4] = "OHiggins"
$row218[0,5] = "Cobresal"
$row218[0,6] = 0
$row218[0,7] = 0
$row218[0,8] = "D"
$row218[0,9] = 3
$row218[0,10] = 3.4
$row218[0,11] = 2.3
$row218[0,12] = 2.1
$row218[0,13] = 3.5
$row218[0,14] = 3.5
$row218[0,15] = -0.25
$row218[0,16] = 1.8
$row218[0,17] = 2.05
$row218[0,18] = 2.75
$row218[0,19] = 1.975
$row218[0,20] = 1.875
$row218[0,21] = -1
$row218[0,22] = 2.5
$row218[0,23] = -1
$row218[0,24] = -0.5
$row218[0,25] = 0.5249999999999999
$row218[0,26] = -1
$row218[0,27] = 0.875
$ws.Range("B218:AC218").Value = $row218

$row227 = New-Object "object[,]" 1,28
$row227[0,0] = 6078263
$row227[0,1] = "Chile Primera Division"
$row227[0,2] = "Chile Primera Division"
$row227[0,3] = [datetime]"2023-12-03T18:00:00"
$row227[0,4] = "Cobresal"
$row227[0,5] = "Universidad de Chile"
$row227[0,6] = 4
$row227[0,7] = 3
$row227[0,8] = "H"
$row227[0,9] = 2.05
$row227[0,10] = 3.4
$row227[0,11] = 3.5
$row227[0,12] = 2.05
$row227[0,13] = 3.6
$row227[0,14] = 3.5
$row227[0,15] = -0.5
$row227[0,16] = 2
$row227[0,17] = 1.8
$row227[0,18] = 2.75
$row227[0,19] = 1.9
$row227[0,20] = 1.9
$row227[0,21] = 1.05
$row227[0,22] = -1
$row227[0,23] = -1
$row227[0,24] = 1
$row227[0,25] = -1
$row227[0,26] = 0.8999999999999999
$row227[0,27] = -1
$ws.Range("B227:AC227").Value = $row227

$row229 = New-Object "object[,]" 1,28
$row229[0,0] = 6077767
$row229[0,1] = "Chile Primera Division"
$row229[0,2] = "Chile Primera Division"
$row229[0,3] = [datetime]"2023-12-03T18:00:00"
$row229[0,4] = "Nublense"
$row229[0,5] = "Huachipato"
$row229[0,6] = 0
$row229[0,7] = 1
$row229[0,8] = "A"
$row229[0,9] = 2.75
$row229[0,10] = 3.4
$row229[0,11] = 2.45
$row229[0,12] = 2.875
$row229[0,13] = 3.3
$row229[0,14] = 2.5
$row229[0,15] = 0
$row229[0,16] = 2.05
$row229[0,17] = 1.8
$row229[0,18] = 2.25
$row229[0,19] = 1.8
$row229[0,20] = 2.05
$row229[0,21] = -1
$row229[0,22] = -1
$row229[0,23] = 1.5
$row229[0,24] = -1
$row229[0,25] = 0.8
$row229[0,26] = -1
$row229[0,27] = 1.05
$ws.Range("B229:AC229").Value = $row229

$row230 = New-Object "object[,]" 1,28
$row230[0,0] = 6078997
$row230[0,1] = "Chile Primera Division"
$row230[0,2] = "Chile Primera Division"
$row230[0,3] = [datetime]"2023-12-08T18:00:00"
$row230[0,4] = "Union Espanola"
$row230[0,5] = "Cobresal"
$row230[0,6] = 1
$row230[0,7] = 0
$row230[0,8] = "H"
$row230[0,9] = 3.8
$row230[0,10] = 3.6
$row230[0,11] = 1.909
$row230[0,12] = 2.7
$row230[0,13] = 3.6
$row230[0,14] = 2.45
$row230[0,15] = 0
$row230[0,16] = 1.975
$row230[0,17] = 1.825
$row230[0,18] = 2.75
$row230[0,19] = 1.775
$row230[0,20] = 2.025
$row230[0,21] = 1.7
$row230[0,22] = -1
$row230[0,23] = -1
$row230[0,24] = 0.9750000000000001
$row230[0,25] = -1
$row230[0,26] = -1
$row230[0,27] = 1.025
$ws.Range("B230:AC230").Value = $row230

$row231 = New-Object "object[,]" 1,28
$row231[0,0] = 6078267
$row231[0,1] = "Chile Primera Division"
$row231[0,2] = "Chile Primera Division"
$row231[0,3] = [datetime]"2023-12-08T18:00:00"
$row231[0,4] = "Huachipato"
$row231[0,5] = "Audax Italiano"
$row231[0,6] = 2
$row231[0,7] = 0
$row231[0,8] = "H"
$row231[0,9] = 1.5
$row231[0,10] = 4.333
$row231[0,11] = 6
$row231[0,12] = 1.444
$row231[0,13] = 4.75
$row231[0,14] = 7
$row231[0,15] = -1.25
$row231[0,16] = 2.025
$row231[0,17] = 1.825
$row231[0,18] = 2.75
$row231[0,19] = 1.8
$row231[0,20] = 2.05
$row231[0,21] = 0.444
$row231[0,22] = -1
$row231[0,23] = -1
$row231[0,24] = 1.025
$row231[0,25] = -1
$row231[0,26] = -1
$row231[0,27] = 1.05
$ws.Range("B231:AC231").Value = $row231

$row232 = New-Object "object[,]" 1,28
$row232[0,0] = 6143704
$row232[0,1] = "Chile Primera Division"
$row232[0,2] = "Chile Primera Division"
$row232[0,3] = [datetime]"2023-12-08T18:00:00"
$row232[0,4] = "Curico Unido"
$row232[0,5] = "Colo Colo"
$row232[0,6] = 0
$row232[0,7] = 1
$row232[0,8] = "A"
$row232[0,9] = 6.5
$row232[0,10] = 4.75
$row232[0,11] = 1.4
$row232[0,12] = 12
$row232[0,13] = 8.5
$row232[0,14] = 1.166
$row232[0,15] = 2
$row232[0,16] = 2
$row232[0,17] = 1.8
$row232[0,18] = 3.25
$row232[0,19] = 1.875
$row232[0,20] = 1.925
$row232[0,21] = -1
$row232[0,22] = -1
$row232[0,23] = 0.1659999999999999
$row232[0,24] = 1
$row232[0,25] = -1
$row232[0,26] = -1
$row232[0,27] = 0.925
$ws.Range("B232:AC232").Value = $row232

$row233 = New-Object "object[,]" 1,28
$row233[0,0] = 6078269
$row233[0,1] = "Chile Primera Division"
$row233[0,2] = "Chile Primera Division"
$row233[0,3] = [datetime]"2023-12-09T18:00:00"
$row233[0,4] = "Universidad de Chile"
$row233[0,5] = "Nublense"
$row233[0,6] = 3
$row233[0,7] = 1
$row233[0,8] = "H"
$row233[0,9] = 1.85
$row233[0,10] = 3.4
$row233[0,11] = 4.333
$row233[0,12] = 1.8
$row233[0,13] = 3.6
$row233[0,14] = 4.5
$row233[0,15] = -0.75
$row233[0,16] = 1.925
$row233[0,17] = 1.925
$row233[0,18] = 2.5
$row233[0,19] = 2.025
$row233[0,20] = 1.825
$row233[0,21] = 0.8
$row233[0,22] = -1
$row233[0,23] = -1
$row233[0,24] = 0.925
$row233[0,25] = -1
$row233[0,26] = 1.025
$row233[0,27] = -1
$ws.Range("B233:AC233").Value = $row233

$row234 = New-Object "object[,]" 1,28
$row234[0,0] = 6077768
$row234[0,1] = "Chile Primera Division"
$row234[0,2] = "Chile Primera Division"
$row234[0,3] = [datetime]"2023-12-09T18:00:00"
$row234[0,4] = "Union La Calera"
$row234[0,5] = "Universidad Catolica"
$row234[0,6] = 0
$row234[0,7] = 3
$row234[0,8] = "A"
$row234[0,9] = 2.05
$row234[0,10] = 3.5
$row234[0,11] = 3.4
$row234[0,12] = 2.05
$row234[0,13] = 3.6
$row234[0,14] = 3.4
$row234[0,15] = -0.25
$row234[0,16] = 1.8
$row234[0,17] = 2
$row234[0,18] = 2.75
$row234[0,19] = 1.975
$row234[0,20] = 1.825
$row234[0,21] = -1
$row234[0,22] = -1
$row234[0,23] = 2.4
$row234[0,24] = -1
$row234[0,25] = 1
$row234[0,26] = 0.4875
$row234[0,27] = -0.5
$ws.Range("B234:AC234").Value = $row234

$row235 = New-Object "object[,]" 1,28
$row235[0,0] = 6078268
$row235[0,1] = "Chile Primera Division"
$row235[0,2] = "Chile Primera Division"
$row235[0,3] = [datetime]"2023-12-09T18:00:00"
$row235[0,4] = "OHiggins"
$row235[0,5] = "Palestino"
$row235[0,6] = 0
$row235[0,7] = 1
$row235[0,8] = "A"
$row235[0,9] = 3.1
$row235[0,10] = 3.3
$row235[0,11] = 2.3
$row235[0,12] = 2.9
$row235[0,13] = 3.4
$row235[0,14] = 2.375
$row235[0,15] = 0.25
$row235[0,16] = 1.8
$row235[0,17] = 2
$row235[0,18] = 2.75
$row235[0,19] = 2
$row235[0,20] = 1.8
$row235[0,21] = -1
$row235[0,22] = -1
$row235[0,23] = 1.375
$row235[0,24] = -1
$row235[0,25] = 1
$row235[0,26] = -1
$row235[0,27] = 0.8
$ws.Range("B235:AC235").Value = $row235

$row236 = New-Object "object[,]" 1,28
$row236[0,0] = 6078998
$row236[0,1] = "Chile Primera Division"
$row236[0,2] = "Chile Primera Division"
$row236[0,3] = [datetime]"2023-12-09T18:00:00"
$row236[0,4] = "Magallanes"
$row236[0,5] = "Coquimbo Unido"
$row236[0,6] = 2
$row236[0,7] = 3
$row236[0,8] = "A"
$row236[0,9] = 1.909
$row236[0,10] = 3.6
$row236[0,11] = 3.8
$row236[0,12] = 2.15
$row236[0,13] = 3.75
$row236[0,14] = 3.1
$row236[0,15] = -0.25
$row236[0,16] = 1.85
$row236[0,17] = 1.95
$row236[0,18] = 3
$row236[0,19] = 1.85
$row236[0,20] = 1.95
$row236[0,21] = -1
$row236[0,22] = -1
$row236[0,23] = 2.1
$row236[0,24] = -1
$row236[0,25] = 0.95
$row236[0,26] = 0.8500000000000001
$row236[0,27] = -1
$ws.Range("B236:AC236").Value = $row236

$row237 = New-Object "object[,]" 1,28
$row237[0,0] = 6077499
$row237[0,1] = "Chile Primera Division"
$row237[0,2] = "Chile Primera Division"
$row237[0,3] = [datetime]"2023-12-09T18:00:00"
$row237[0,4] = "Deportes Copiapo"
$row237[0,5] = "Everton de Vina"
$row237[0,6] = 2
$row237[0,7] = 0
$row237[0,8] = "H"
$row237[0,9] = 2.1
$row237[0,10] = 3.4
$row237[0,11] = 3.4
$row237[0,12] = 2.9
$row237[0,13] = 3.4
$row237[0,14] = 2.4
$row237[0,15] = 0.25
$row237[0,16] = 1.775
$row237[0,17] = 2.1
$row237[0,18] = 2.75
$row237[0,19] = 1.85
$row237[0,20] = 2
$row237[0,21] = 1.9
$row237[0,22] = -1
$row237[0,23] = -1
$row237[0,24] = 0.7749999999999999
$row237[0,25] = -1
$row237[0,26] = -1
$row237[0,27] = 1
$ws.Range("B237:AC237").Value = $row237

# --- Partial cell updates (upcoming-fixture rows with refreshed odds / swapped fixtures) ---
$ws.Range("N238").Value = 2.5
$ws.Range("O238").Value = 3.25

$ws.Range("N239").Value = 2.3
$ws.Range("R239").Value = 1.9
$ws.Range("S239").Value = 1.95
$ws.Range("U239").Value = 2.05
$ws.Range("V239").Value = 1.8

$ws.Range("N240").Value = 4.333
$ws.Range("O240").Value = 3.75
$ws.Range("P240").Value = 1.8
$ws.Range("Q240").Value = 0.75
$ws.Range("R240").Value = 1.825
$ws.Range("S240").Value = 2.025
$ws.Range("U240").Value = 1.875
$ws.Range("V240").Value = 1.975

$ws.Range("B241").Value = 7723513
$ws.Range("F241").Value = "Union La Calera"
$ws.Range("G241").Value = "Universidad Catolica"
$ws.Range("K241").Value = 2.75
$ws.Range("L241").Value = 3.25
$ws.Range("M241").Value = 2.6
$ws.Range("N241").Value = 2.4
$ws.Range("O241").Value = 3.25
$ws.Range("P241").Value = 3.1
$ws.Range("Q241").Value = -0.25
$ws.Range("R241").Value = 2.05
$ws.Range("S241").Value = 1.8
$ws.Range("T241").Value = 2.5
$ws.Range("U241").Value = 1.975
$ws.Range("V241").Value = 1.875

$ws.Range("B242").Value = 7723518
$ws.Range("F242").Value = "Deportes Copiapo"
$ws.Range("G242").Value = "OHiggins"
$ws.Range("K242").Value = 2.5
$ws.Range("L242").Value = 3.2
$ws.Range("M242").Value = 2.875
$ws.Range("N242").Value = 2.45
$ws.Range("Q242").Value = -0.25
$ws.Range("R242").Value = 2.125
$ws.Range("S242").Value = 1.75
$ws.Range("U242").Value = 2.05
$ws.Range("V242").Value = 1.8

$ws.Range("N243").Value = 1.85
$ws.Range("P243").Value = 4.333

$ws.Range("P244").Value = 2.875
$ws.Range("R244").Value = 1.875
$ws.Range("S244").Value = 1.975
$ws.Range("U244").Value = 2
$ws.Range("V244").Value = 1.85

Write-Output "edit complete"